# Update schedule values for several training trials.
# Columns: A trialTrain, B x_fixStart, C y_fixStart, D x_corrSteps,
#          E y_corrSteps, F x_nrSteps, G y_nrSteps, H alienID, I praclen, J version
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of sheet row number -> new values for columns E (y_corrSteps), G (y_nrSteps), H (alienID)
$updates = @{
    4  = @{ E = 6; G = 3; H = 13 }
    8  = @{ E = 6; G = 3; H = 13 }
    16 = @{ E = 7; G = 3; H = 13 }
    18 = @{ E = 6; G = 3; H = 13 }
    23 = @{ E = 5; G = 3; H = 13 }
    27 = @{ E = 7; G = 3; H = 13 }
}

foreach ($r in $updates.Keys) {
    $vals = $updates[$r]
    $ws.Cells.Item($r, 5).Value = $vals.E
    $ws.Cells.Item($r, 7).Value = $vals.G
    $ws.Cells.Item($r, 8).Value = $vals.H
}
